# Updates cryptocurrency Price (D) and Volume(1h) (E) columns for rows
# 2-51, matching the latest scrape from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of Price values are plain decimals (e.g. "235.14"). Left alone,
# assigning a numeric-looking string to .Value makes Excel auto-convert the
# cell to a Number (dropping e.g. a trailing zero: "0.0990" -> 0.099). Mark
# just those specific cells as Text first so the literal source string is
# kept, consistent with the other "thousands.dot" price cells that already
# stay text automatically because they are not valid numbers.
$textCells = @(
    "D5", "D7", "D13", "D14", "D16", "D19", "D20", "D22", "D24", "D26",
    "D30", "D34", "D35", "D36", "D39", "D40", "D41", "D46", "D47", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '37.856.41'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '2.086.35'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '235.14'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").Value = '59.53'
$ws.Range("E7").Value = '  +3.52%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -1.06%  '
$ws.Range("E10").Value = '  +2.55%  '
$ws.Range("E11").Value = '  +3.05%  '
$ws.Range("D12").Value = '2.392.97'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").Value = '14.74'
$ws.Range("E13").Value = '  +1.94%  '
$ws.Range("D14").Value = '21.40'
$ws.Range("E14").Value = '  +3.12%  '
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").Value = '5.32'
$ws.Range("E16").Value = '  +2.38%  '
$ws.Range("D17").Value = '2.090.50'
$ws.Range("E17").Value = '  +0.96%  '
$ws.Range("D18").Value = '37.779.78'
$ws.Range("E18").Value = '  +1.14%  '
$ws.Range("D19").Value = '6.23'
$ws.Range("E19").Value = '  -2.16%  '
$ws.Range("D20").Value = '71.79'
$ws.Range("E20").Value = '  +2.90%  '
$ws.Range("E21").Value = '  +1.48%  '
$ws.Range("D22").Value = '229.21'
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '2.42'
$ws.Range("E24").Value = '  -0.58%  '
$ws.Range("E25").Value = '  +0.81%  '
$ws.Range("D26").Value = '170.76'
$ws.Range("E26").Value = '  +2.16%  '
$ws.Range("E27").Value = '  +9.52%  '
$ws.Range("E28").Value = '  +2.57%  '
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("D30").Value = '19.59'
$ws.Range("E30").Value = '  +2.44%  '
$ws.Range("E31").Value = '  +1.57%  '
$ws.Range("E32").Value = '  +4.46%  '
$ws.Range("E33").Value = '  +2.26%  '
$ws.Range("D34").Value = '4.70'
$ws.Range("E34").Value = '  +2.87%  '
$ws.Range("D35").Value = '2.52'
$ws.Range("E35").Value = '  +0.81%  '
$ws.Range("D36").Value = '3.52'
$ws.Range("E36").Value = '  +6.54%  '
$ws.Range("E37").Value = '  +2.50%  '
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("D39").Value = '5.45'
$ws.Range("E39").Value = '  -3.83%  '
$ws.Range("D40").Value = '0.0990'
$ws.Range("E40").Value = '  +2.44%  '
$ws.Range("D41").Value = '99.82'
$ws.Range("E41").Value = '  +1.53%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  +1.15%  '
$ws.Range("D44").Value = '1.464.03'
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("E45").Value = '  +0.66%  '
$ws.Range("D46").Value = '4.31'
$ws.Range("E46").Value = '  +4.95%  '
$ws.Range("D47").Value = '16.18'
$ws.Range("E47").Value = '  +5.49%  '
$ws.Range("E48").Value = '  +3.93%  '
$ws.Range("E49").Value = '  +3.19%  '
$ws.Range("E50").Value = '  +2.70%  '
$ws.Range("D51").Value = '47.70'
$ws.Range("E51").Value = '  +6.37%  '

# Drop the temporary Text format again so the cells keep the workbook's
# default (unstyled) look, same as before the update.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

